# Swap the values of columns B, D, E, F, G between each of the paired rows.
# These pairs correspond to duplicate stock-item lines whose batch figures
# were transposed between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(127, 128),
    @(227, 228),
    @(232, 233),
    @(243, 244),
    @(322, 323),
    @(380, 381),
    @(382, 383),
    @(385, 386),
    @(442, 443),
    @(473, 474)
)

$cols = @("B", "D", "E", "F", "G")

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"

        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2

        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}
